$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 504, shifting the existing
# rows 504:533 down to 506:535 (matches the diff: dimension grows from
# R533 to R535, and two brand-new "Crespo record" observations are added
# right after row 503).
$ws.Rows.Item(504).Insert()
$ws.Rows.Item(504).Insert()

# New row 504 - "Primera" quality observation dated 45223
$ws.Cells.Item(504, 1).Value = 7
$ws.Cells.Item(504, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(504, 3).Value = "Ñuble"
$ws.Cells.Item(504, 4).Value = 45223
$ws.Cells.Item(504, 5).Value = 16
$ws.Cells.Item(504, 6).Value = 100112006
$ws.Cells.Item(504, 7).Value = "Repollo"
$ws.Cells.Item(504, 8).Value = "Crespo record"
$ws.Cells.Item(504, 9).Value = "Primera"
$ws.Cells.Item(504, 10).Value = 600
$ws.Cells.Item(504, 11).Value = 1200
$ws.Cells.Item(504, 12).Value = 1400
$ws.Cells.Item(504, 13).Value = 1300
$ws.Cells.Item(504, 14).Value = "`$/unidad"
$ws.Cells.Item(504, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(504, 16).Value = 1300
$ws.Cells.Item(504, 17).Value = 1
$ws.Cells.Item(504, 18).Value = "Hortaliza"

# New row 505 - "Segunda" quality observation, same date
$ws.Cells.Item(505, 1).Value = 7
$ws.Cells.Item(505, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(505, 3).Value = "Ñuble"
$ws.Cells.Item(505, 4).Value = 45223
$ws.Cells.Item(505, 5).Value = 16
$ws.Cells.Item(505, 6).Value = 100112006
$ws.Cells.Item(505, 7).Value = "Repollo"
$ws.Cells.Item(505, 8).Value = "Crespo record"
$ws.Cells.Item(505, 9).Value = "Segunda"
$ws.Cells.Item(505, 10).Value = 500
$ws.Cells.Item(505, 11).Value = 1000
$ws.Cells.Item(505, 12).Value = 1000
$ws.Cells.Item(505, 13).Value = 1000
$ws.Cells.Item(505, 14).Value = "`$/unidad"
$ws.Cells.Item(505, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(505, 16).Value = 1000
$ws.Cells.Item(505, 17).Value = 1
$ws.Cells.Item(505, 18).Value = "Hortaliza"
